$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = '26.448.89'
$ws.Range("D3").Value = '1.626.52'
$ws.Range("E3").Value = '  -0.59%  '
$ws.Range("E4").Value = '  +0.19%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '213.11'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.03%  '
$ws.Range("E6").Value = '  +1.56%  '
$ws.Range("E7").Value = '  +0.18%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.0623'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.27%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.248'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.31%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '18.88'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.02%  '
$ws.Range("E11").Value = '  +1.15%  '
$ws.Range("D12").Value = '1.854.04'
$ws.Range("E12").Value = '  -0.57%  '
$ws.Range("D13").Value = '1.625.52'
$ws.Range("E13").Value = '  -1.49%  '
$ws.Range("E14").Value = '  +1.79%  '
$ws.Range("E15").Value = '  -0.20%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.88'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.97%  '
$ws.Range("D17").Value = '26.512.73'
$ws.Range("E17").Value = '  -0.48%  '
$ws.Range("D18").Value = '0.0₃0739'
$ws.Range("E18").Value = '  -0.09%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '215.19'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.45%  '
$ws.Range("E20").Value = '  +0.17%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.26'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.11%  '
$ws.Range("E23").Value = '  -1.12%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.01'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +4.85%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '148.54'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.72%  '
$ws.Range("E26").Value = '  +0.13%  '
$ws.Range("E27").Value = '  +0.10%  '
$ws.Range("E28").Value = '  +1.86%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.55'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.01%  '
$ws.Range("E30").Value = '  -1.54%  '
$ws.Range("E31").Value = '  -0.68%  '
$ws.Range("E32").Value = '  +2.97%  '
$ws.Range("E33").Value = '  -0.29%  '
$ws.Range("E34").Value = '  -0.85%  '
$ws.Range("E35").Value = '  -1.10%  '
$ws.Range("D36").Value = '1.219.04'
$ws.Range("E36").Value = '  +4.34%  '
$ws.Range("E37").Value = '  +4.02%  '
$ws.Range("E38").Value = '  +0.17%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.794'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.92%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.505'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.82%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.27'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.30%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.791'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.36%  '
$ws.Range("E43").Value = '  +0.02%  '
$ws.Range("D44").Value = '1.765.00'
$ws.Range("E44").Value = '  -0.49%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '93.12'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.00%  '
$ws.Range("E46").Value = '  +1.53%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '54.86'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.44%  '
$ws.Range("E48").Value = '  -0.60%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.54'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.13%  '
$ws.Range("E50").Value = '  -0.78%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.01'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.39%  '
